# Release MHD 4.2.2 close #419
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version bump
$ws.Range("B3").Value = "4.2.2"

# Publication date bump
$ws.Range("B8").Value = "2024-05-18T12:39:23-05:00"

# Contact rows now carry real contact detail display values
$ws.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$ws.Range("B11").Value = "null (iti@ihe.net)"
$ws.Range("B12").Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"
